$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from H1 into the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 6
